$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.216.99"
$ws.Range("E2").Value = "  +3.14%  "
$ws.Range("D3").Value = "'3.148.87"
$ws.Range("E3").Value = "  +2.22%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'591.66"
$ws.Range("E5").Value = "  +1.55%  "
$ws.Range("D6").Value = "'147.66"
$ws.Range("E6").Value = "  +2.76%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'3.142.26"
$ws.Range("E8").Value = "  +2.30%  "
$ws.Range("D9").Value = "'0.535"
$ws.Range("E9").Value = "  +0.90%  "
$ws.Range("E10").Value = "  +14.55%  "
$ws.Range("D11").Value = "'5.74"
$ws.Range("E11").Value = "  -0.21%  "
$ws.Range("D12").Value = "'0.471"
$ws.Range("E12").Value = "  +0.86%  "
$ws.Range("D13").Value = "'0.0000255"
$ws.Range("E13").Value = "  +5.39%  "
$ws.Range("D14").Value = "'37.47"
$ws.Range("E14").Value = "  +5.76%  "
$ws.Range("E15").Value = "  -0.74%  "
$ws.Range("D16").Value = "'3.664.26"
$ws.Range("E16").Value = "  +2.22%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "'64.028.82"
$ws.Range("E17").Value = "  +2.98%  "
$ws.Range("B18").Value = "Polkadot"
$ws.Range("C18").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D18").Value = "'7.20"
$ws.Range("E18").Value = "  -1.33%  "
$ws.Range("D19").Value = "'3.140.84"
$ws.Range("E19").Value = "  +2.14%  "
$ws.Range("D20").Value = "'468.59"
$ws.Range("E20").Value = "  +4.24%  "
$ws.Range("D21").Value = "'14.40"
$ws.Range("E21").Value = "  +2.58%  "
$ws.Range("D22").Value = "'0.739"
$ws.Range("E22").Value = "  +0.73%  "
$ws.Range("D23").Value = "'7.59"
$ws.Range("E23").Value = "  +1.58%  "
$ws.Range("D24").Value = "'13.39"
$ws.Range("E24").Value = "  -2.63%  "
$ws.Range("D25").Value = "'82.52"
$ws.Range("E25").Value = "  +0.81%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").Value = "'9.01"
$ws.Range("E27").Value = "  +8.68%  "
$ws.Range("E28").Value = "  +1.76%  "
$ws.Range("E29").Value = "  -0.76%  "
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("E31").Value = "  +1.09%  "
$ws.Range("D32").Value = "'27.29"
$ws.Range("E33").Value = "  -3.23%  "
$ws.Range("D34").Value = "'0.0₃0885"
$ws.Range("E34").Value = "  +11.10%  "
$ws.Range("D35").Value = "'2.40"
$ws.Range("E35").Value = "  +8.45%  "
$ws.Range("B36").Value = "Mantle"
$ws.Range("C36").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D36").Value = "'1.06"
$ws.Range("E36").Value = "  +2.07%  "
$ws.Range("B37").Value = "dogwifhat"
$ws.Range("C37").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D37").Value = "'3.45"
$ws.Range("E37").Value = "  +16.02%  "
$ws.Range("E38").Value = "  +1.57%  "
$ws.Range("D39").Value = "'466.54"
$ws.Range("E39").Value = "  +9.74%  "
$ws.Range("D40").Value = "'51.03"
$ws.Range("E40").Value = "  +1.17%  "
$ws.Range("D41").Value = "'8.78"
$ws.Range("E41").Value = "  -0.21%  "
$ws.Range("E42").Value = "  +1.71%  "
$ws.Range("D43").Value = "'2.916.27"
$ws.Range("E43").Value = "  +0.22%  "
$ws.Range("E44").Value = "  +0.64%  "
$ws.Range("E46").Value = "  +2.01%  "
$ws.Range("D47").Value = "'126.95"
$ws.Range("E47").Value = "  +1.97%  "
$ws.Range("D48").Value = "'35.62"
$ws.Range("E48").Value = "  +2.45%  "
$ws.Range("E49").Value = "  +0.03%  "
$ws.Range("E50").Value = "  +0.67%  "
$ws.Range("D51").Value = "'24.90"
$ws.Range("E51").Value = "  +1.22%  "
